# Weekly update: two new price records (row 15 and row 16) are added at the
# top of the data table, pushing every existing record down by two rows.
# The last two existing records end up appended as brand-new rows 98 and 99.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make sure the two brand-new rows at the bottom (98 and 99) get the same
# date number format used throughout column D before we copy values into
# them (setting the format first avoids Excel auto-assigning a generic
# date format and bloating the style table).
$ws.Range("D98").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D99").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Shift every data row down by two, starting from the bottom so that a
# source row is always read before it gets overwritten.
for ($r = 99; $r -ge 17; $r--) {
    $srcRow = $r - 2
    $src = $ws.Range("A" + $srcRow + ":R" + $srcRow).Value()
    $ws.Range("A" + $r + ":R" + $r).Value = $src
}

# Row 15: brand new record for 2023-05-25 (serial 45061).
$ws.Range("A15").Value = 10
$ws.Range("B15").Value = "Vega Modelo de Temuco"
$ws.Range("C15").Value = "La Araucanía"
$ws.Range("D15").Value = 45061
$ws.Range("E15").Value = 9
$ws.Range("F15").Value = 100112010
$ws.Range("G15").Value = "Achicoria"
$ws.Range("H15").Value = "Sin especificar"
$ws.Range("I15").Value = "Primera"
$ws.Range("J15").Value = 65
$ws.Range("K15").Value = 10000
$ws.Range("L15").Value = 10000
$ws.Range("M15").Value = 10000
$ws.Range("N15").Value = "`$/caja 18 unidades"
$ws.Range("O15").Value = "Región Metropolitana"
$ws.Range("P15").Value = 556
$ws.Range("Q15").Value = 18
$ws.Range("R15").Value = "Hortaliza"

# Row 16: brand new record also for 2023-05-25 (serial 45061).
$ws.Range("A16").Value = 10
$ws.Range("B16").Value = "Vega Modelo de Temuco"
$ws.Range("C16").Value = "La Araucanía"
$ws.Range("D16").Value = 45061
$ws.Range("E16").Value = 9
$ws.Range("F16").Value = 100112010
$ws.Range("G16").Value = "Achicoria"
$ws.Range("H16").Value = "Sin especificar"
$ws.Range("I16").Value = "Primera"
$ws.Range("J16").Value = 90
$ws.Range("K16").Value = 7000
$ws.Range("L16").Value = 7000
$ws.Range("M16").Value = 7000
$ws.Range("N16").Value = "`$/caja 18 unidades"
$ws.Range("O16").Value = "Región del Maule"
$ws.Range("P16").Value = 389
$ws.Range("Q16").Value = 18
$ws.Range("R16").Value = "Hortaliza"
